{"js": "// Append two new paragraphs to the end of the document body:\n//   1) \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\"  (date line, split the same way the other\n//      date paragraphs in this diary are split: leading \"2022\" run, then\n//      the CJK remainder as its own run)\n//   2) \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002\" (weather/body line)\n//\n// We build the new paragraphs from an explicit OOXML fragment (via\n// Range.insertOoxml) rather than Body.insertParagraph/insertText so the\n// inserted runs/paragraph-mark formatting land exactly as authored\n// (matching the surrounding paragraphs' \"<w:rFonts w:hint=\"eastAsia\"/>\"\n// convention on the CJK runs) instead of being auto-merged/re-hinted.\n\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange(\"Whole\");\n\nconst flatOpcNamespace = \"http://schemas.microsoft.com/office/2006/xmlPackage\";\n\nconst newParagraphsXml =\n  \"<w:p>\" +\n    \"<w:r><w:t>2022</w:t></w:r>\" +\n    \"<w:r><w:rPr><w:rFonts w:hint=\\\"eastAsia\\\"/></w:rPr><w:t>\u5e746\u67087\u65e5\u661f\u671f\u4e8c</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p>\" +\n    \"<w:pPr><w:rPr><w:rFonts w:hint=\\\"eastAsia\\\"/></w:rPr></w:pPr>\" +\n    \"<w:r><w:rPr><w:rFonts w:hint=\\\"eastAsia\\\"/></w:rPr><w:t>\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002</w:t></w:r>\" +\n  \"</w:p>\";\n\nconst flatOpcXml =\n  \"<?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?>\" +\n  \"<pkg:package xmlns:pkg=\\\"\" + flatOpcNamespace + \"\\\">\" +\n    \"<pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\">\" +\n      \"<pkg:xmlData>\" +\n        \"<w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\">\" +\n          \"<w:body>\" + newParagraphsXml + \"</w:body>\" +\n        \"</w:document>\" +\n      \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ninsertionRange.insertOoxml(flatOpcXml, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Append two new paragraphs to the end of the document body:\n#   1) \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\"  (date line, split the same way the other\n#      date paragraphs in this diary are split: leading \"2022\" run, then\n#      the CJK remainder as its own run)\n#   2) \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002\" (weather/body line)\n#\n# Built from an explicit OOXML fragment (Range.InsertXML with a minimal\n# flat-OPC wrapper) so the new runs / paragraph-mark formatting land\n# exactly as authored (matching the surrounding paragraphs' eastAsia\n# font-hint convention) rather than being merged/re-hinted by plain\n# TypeText/InsertAfter calls.\n\n$d = $word.ActiveDocument\n\n$newParagraphsXml = '<w:p><w:r><w:t>2022</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u5e746\u67087\u65e5\u661f\u671f\u4e8c</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002</w:t></w:r></w:p>'\n\n$flatOpcXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# A Range built fresh from Document.Range(start, end) at the very end of\n# the story (rather than a paragraph's Range collapsed via Collapse())\n# inserts after that point without touching/replacing the existing last\n# paragraph's text or its paragraph mark.\n$endPos = $d.Content.End\n$insertionRange = $d.Range($endPos, $endPos)\n$insertionRange.InsertXML($flatOpcXml)\n"}
